$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Delete the empty trailing slide (sldId 259 / slide4.xml).
#    It has no shapes - it's a blank placeholder slide that was removed.
# ---------------------------------------------------------------------
for ($i = $p.Slides.Count; $i -ge 1; $i--) {
    $slide = $p.Slides.Item($i)
    if (($slide.SlideID -eq 259) -or ($slide.Shapes.Count -eq 0)) {
        $slide.Delete()
    }
}

# ---------------------------------------------------------------------
# 2. Refresh the cached "datetimeFigureOut" date text (8/27/2020 ->
#    8/28/2020) on the slide master and on every slide layout's Date
#    placeholder.
# ---------------------------------------------------------------------
$newDate = "8/28/2020"
$oldDate = "8/27/2020"
$ppPlaceholderDate = 16

$master = $p.SlideMaster

for ($j = 1; $j -le $master.Shapes.Count; $j++) {
    $shape = $master.Shapes.Item($j)
    if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
        if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
            $shape.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($i = 1; $i -le $master.CustomLayouts.Count; $i++) {
    $layout = $master.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $shape = $layout.Shapes.Item($j)
        if ($shape.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
            if ($shape.TextFrame.TextRange.Text -eq $oldDate) {
                $shape.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
